$wb = $excel.ActiveWorkbook

# Fix typo in sheet name: "SACRU - URLs" -> "SARCU - URLs"
$wsSarcu = $wb.Worksheets.Item("SACRU - URLs")
$wsSarcu.Name = "SARCU - URLs"

# Re-select cell A22 on the renamed sheet
$wsSarcu.Activate()
$wsSarcu.Range("A22").Select()

# Size column A to fit the long URL on the "FlyingPartnerAir URL" sheet
$wsPartnerUrl = $wb.Worksheets.Item("OUATT - FlyingPartnerAir URL")
$wsPartnerUrl.Activate()
$wsPartnerUrl.Columns("A:A").AutoFit()
$wsPartnerUrl.Columns("A:A").ColumnWidth = 129.3

# Visit "FAQ Header Names" (keeps its own cached G3 selection, but is no longer the active tab)
$wsFaq = $wb.Worksheets.Item("FAQ Header Names")
$wsFaq.Activate()

# Make "TUC Dropdown List Count" the final active sheet/tab (0-based activeTab = 1)
$wsCount = $wb.Worksheets.Item("TUC Dropdown List Count")
$wsCount.Activate()
